$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in row 10 (Senado 18) color: #348cae4 -> #348cae
$ws.Range("C10").Value = "#348cae"

# Insert a new row before current row 15 (Distrito local 21) to host the
# new "Gobernatura 21" entry, shifting the remaining rows down.
$ws.Rows("15").Insert()

$ws.Range("A15").Value = "Gobernatura 21"
$ws.Range("B15").Value = "gb_21"
$ws.Range("C15").Value = "#588157"

# Update the selection to match the target state
$ws.Range("C15").Select()
